# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colour scheme (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral" colour scheme (used by the slide master / all slides)
# The authored change swaps the two themes' contents: theme1.xml becomes the
# "Integral" theme and theme2.xml becomes the "Office Theme" theme (file names /
# relationships stay put - only the colour definitions move).
#
# The slide master's theme (reached through Designs(1).SlideMaster.Theme) is the
# theme actually driving the visible deck, and it is the one exposed for editing
# by the PowerPoint object model here, so we repaint its 12-slot theme colour
# scheme with the "Office Theme" colours that the target theme2.xml should hold.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colours (target content for ppt/theme/theme2.xml), in the
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order that ThemeColorScheme
# exposes as Item(1..12).
$colors.Item(1).RGB  = 0       # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1    FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2    44546A
$colors.Item(4).RGB  = 15132391  # lt2    E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1 5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2 ED7D31
$colors.Item(7).RGB  = 10855845  # accent3 A5A5A5
$colors.Item(8).RGB  = 49407     # accent4 FFC000
$colors.Item(9).RGB  = 12874308  # accent5 4472C4
$colors.Item(10).RGB = 4697456   # accent6 70AD47
$colors.Item(11).RGB = 12673797  # hlink   0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
